$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 44134
$ws.Range("C2").Value = 402401
$ws.Range("D2").Value = 9711
$ws.Range("E2").Value = 48655
$ws.Range("F2").Value = 2327
$ws.Range("G2").Value = 12.09
$ws.Range("H2").Value = 23.96

# Row 3
$ws.Range("C3").Value = 182270
$ws.Range("D3").Value = 5705

# Row 4
$ws.Range("B4").Value = 44134
$ws.Range("C4").Value = 106573
$ws.Range("D4").Value = 2366
$ws.Range("E4").Value = 4104
$ws.Range("H4").Value = 3.04
$ws.Range("K4").Value = 70988
$ws.Range("L4").Value = 2334

# Row 5
$ws.Range("B5").Value = 44134
$ws.Range("C5").Value = 271830
$ws.Range("D5").Value = 4332
$ws.Range("E5").Value = 49715
$ws.Range("F5").Value = 1218
$ws.Range("G5").Value = 22.64
$ws.Range("H5").Value = 29.53
$ws.Range("K5").Value = 219550
$ws.Range("L5").Value = 4124

# Row 7
$ws.Range("B7").NumberFormat = "YYYY-MM-DD"
$ws.Range("B7").Value = 44133
$ws.Range("C7").Value = 505431
$ws.Range("D7").Value = 25804
$ws.Range("F7").Value = 6474
$ws.Range("H7").Value = 26.64
$ws.Range("L7").Value = 24300
$ws.Range("O7").Value = "Success!"

# Row 8
$ws.Range("C8").Formula = "'255208"
$ws.Range("D8").Formula = "'19340"

# Row 14
$ws.Range("B14").Value = 44134
$ws.Range("C14").Value = 259488
$ws.Range("D14").Value = 3341
$ws.Range("E14").Value = 42808
$ws.Range("F14").Value = 798
$ws.Range("G14").Value = 16.5
$ws.Range("H14").Value = 23.89

# Row 16
$ws.Range("B16").Value = 44134
$ws.Range("C16").Value = 144314
$ws.Range("D16").Value = 3990
$ws.Range("E16").Value = 44529
$ws.Range("F16").Value = 1618
$ws.Range("G16").Value = 36.23
$ws.Range("H16").Value = 40.7
$ws.Range("K16").Value = 122896
$ws.Range("L16").Value = 3975

# Row 17
$ws.Range("B17").Value = 44134
$ws.Range("C17").Value = 162720
$ws.Range("D17").Value = 2735
$ws.Range("E17").Value = 37962
$ws.Range("F17").Value = 1035
$ws.Range("G17").Value = 33.57
$ws.Range("H17").Value = 39
$ws.Range("K17").Value = 113089
$ws.Range("L17").Value = 2654

# Row 18
$ws.Range("B18").Value = 44134
$ws.Range("C18").Value = 110874
$ws.Range("D18").Value = 1900
$ws.Range("E18").Value = 22307
$ws.Range("F18").Value = 346
$ws.Range("G18").Value = 22.34
$ws.Range("H18").Value = 19.7
$ws.Range("K18").Value = 99869
$ws.Range("L18").Value = 1756

# Row 20
$ws.Range("B20").Value = 44133
$ws.Range("C20").Value = 67519
$ws.Range("D20").Value = 3708
$ws.Range("E20").Value = 8329
$ws.Range("G20").Value = 12.34
$ws.Range("H20").Value = 18.18

# Row 21
$ws.Range("B21").Value = 44134
$ws.Range("C21").Value = 174303
$ws.Range("D21").Value = 7239
$ws.Range("E21").Value = 29679
$ws.Range("F21").Value = 2598
$ws.Range("G21").Value = 17.03
$ws.Range("H21").Value = 35.89

# Row 23
$ws.Range("B23").Formula = "'2020-10-30"
$ws.Range("C23").Value = 121495
$ws.Range("D23").Value = 1326
$ws.Range("E23").Value = 8273.809499999999
$ws.Range("F23").Value = 79.95780000000001
$ws.Range("G23").Value = 6.81
$ws.Range("H23").Value = 6.03
$ws.Range("K23").Value = 98945.52800000001
$ws.Range("L23").Value = 1228.0086

# Row 24
$ws.Range("B24").Value = 44134
$ws.Range("C24").Value = 31916
$ws.Range("D24").Value = 364
$ws.Range("E24").Value = 114

# Row 26
$ws.Range("B26").Value = 44134
$ws.Range("C26").Value = 11314
$ws.Range("D26").Value = 77
$ws.Range("E26").Value = 511
$ws.Range("G26").Value = 4.66
$ws.Range("H26").Value = 3.9
$ws.Range("K26").Value = 10972
$ws.Range("L26").Value = 77

# Row 28
$ws.Range("B28").Value = 44134
$ws.Range("C28").Value = 157146
$ws.Range("D28").Value = 9975
$ws.Range("E28").Value = 13834
$ws.Range("F28").Value = 788
$ws.Range("G28").Value = 8.800000000000001
$ws.Range("H28").Value = 7.9

# Row 29
$ws.Range("B29").Value = 44134
$ws.Range("C29").Value = 185552
$ws.Range("D29").Value = 3643
$ws.Range("E29").Value = 15993
$ws.Range("G29").Value = 13.13
$ws.Range("K29").Value = 121765

# Row 30
$ws.Range("B30").Value = 44134
$ws.Range("C30").Value = 46424
$ws.Range("E30").Value = 7342
$ws.Range("G30").Value = 21.49
$ws.Range("K30").Value = 34163

# Row 31
$ws.Range("C31").Value = 125166
$ws.Range("D31").Value = 1706
$ws.Range("E31").Value = 5724
$ws.Range("G31").Value = 4.57

# Row 33
$ws.Range("B33").Value = 44134
$ws.Range("C33").Formula = "'112932"
$ws.Range("D33").Formula = "'601"
$ws.Range("E33").Formula = "'1928"
$ws.Range("H33").Value = 1.33

# Row 38
$ws.Range("B38").Value = 44134
$ws.Range("C38").Value = 57928
$ws.Range("D38").Value = 17934
$ws.Range("E38").Value = 9533
$ws.Range("F38").Value = 1949
$ws.Range("G38").Value = 0.16

# Row 41
$ws.Range("B41").Value = 44133
$ws.Range("C41").Value = 916918
$ws.Range("D41").Value = 17571
$ws.Range("E41").Value = 27387
$ws.Range("F41").Value = 1296
$ws.Range("K41").Value = 650369
$ws.Range("L41").Value = 17393

# Row 43
$ws.Range("B43").Value = 44134
$ws.Range("C43").Value = 2155
$ws.Range("G43").Value = 9.859999999999999
$ws.Range("K43").Value = 2059

# Row 44
$ws.Range("B44").Value = 44134
$ws.Range("C44").Value = 45909
$ws.Range("D44").Value = 1007
$ws.Range("E44").Value = 765
$ws.Range("G44").Value = 1.67

# Row 45
$ws.Range("B45").Value = 44134
$ws.Range("C45").Value = 789714
$ws.Range("D45").Value = 16720
$ws.Range("E45").Value = 103732
$ws.Range("F45").Value = 2893
$ws.Range("G45").Value = 13.14
$ws.Range("H45").Value = 17.3

# Row 46
$ws.Range("B46").Value = 44134
$ws.Range("C46").Value = 17144
$ws.Range("D46").Value = 646
$ws.Range("E46").Value = 8479
$ws.Range("F46").Value = 484
$ws.Range("G46").Value = 49.46
$ws.Range("H46").Value = 74.92

# Row 48
$ws.Range("B48").Value = 44134
$ws.Range("C48").Value = 104426
$ws.Range("D48").Value = 2278
$ws.Range("E48").Value = 3676
$ws.Range("G48").Value = 4.29
$ws.Range("H48").Value = 6.05
$ws.Range("K48").Value = 85765
$ws.Range("L48").Value = 2231

# Row 49
$ws.Range("B49").Value = 44134
$ws.Range("C49").Value = 56369
$ws.Range("D49").Value = 888
$ws.Range("E49").Value = 1987
$ws.Range("H49").Value = 4.07
$ws.Range("K49").Value = 46790
$ws.Range("L49").Value = 885

# Row 52
$ws.Range("B52").Value = 44134
$ws.Range("C52").Value = 244045
$ws.Range("D52").Value = 5934
$ws.Range("E52").Value = 7497
$ws.Range("H52").Value = 3.48
$ws.Range("K52").Value = 171678
$ws.Range("L52").Value = 5289

# Row 54
$ws.Range("B54").Value = 44134
$ws.Range("C54").Value = 179612
$ws.Range("D54").Value = 3643
$ws.Range("E54").Value = 35107
$ws.Range("F54").Value = 972
$ws.Range("G54").Value = 19.55

# Row 56
$ws.Range("B56").Value = 44134
$ws.Range("C56").Value = 175893
$ws.Range("D56").Value = 4050
$ws.Range("E56").Value = 14398
$ws.Range("G56").Value = 8.19
$ws.Range("H56").Value = 11.78

# Row 57
$ws.Range("B57").Value = 44134
$ws.Range("C57").Value = 10884
$ws.Range("E57").Value = 416
$ws.Range("G57").Value = 4.43
$ws.Range("K57").Value = 9393

# Row 58
$ws.Range("B58").Value = 44134
$ws.Range("C58").Value = 105242
$ws.Range("D58").Value = 1476
$ws.Range("E58").Value = 9029
$ws.Range("F58").Value = 162
$ws.Range("G58").Value = 11.01
$ws.Range("H58").Value = 11.86
$ws.Range("K58").Value = 82015
$ws.Range("L58").Value = 1374

# Row 59
$ws.Range("B59").Value = 44133
$ws.Range("C59").Value = 306327
$ws.Range("D59").Value = 7056
$ws.Range("E59").Value = 8886
$ws.Range("F59").Value = 637
$ws.Range("G59").Value = 4.78
$ws.Range("H59").Value = 9.58
$ws.Range("K59").Value = 185731
$ws.Range("L59").Value = 6649
